$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: the match-status enumeration in the "Business Rules" section —
# "aktív(vagy folyamatban, [ongoing])" becomes simply "elkezdett".
# ---------------------------------------------------------------------
$rStatus = $d.Content
$null = $rStatus.Find.Execute(
    "“aktív(vagy folyamatban, [ongoing])”",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "“elkezdett”",
    2)

# ---------------------------------------------------------------------
# Edit 2: the "JÁTÉKOS:" heading is renamed to "FELHASZNÁLÓ:" — authored
# as two runs (label text + trailing colon), both bold, matching the
# canonical OOXML produced by the original edit.
# ---------------------------------------------------------------------
$rHeading = $d.Content
$null = $rHeading.Find.Execute(
    "JÁTÉKOS:",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "FELHASZNÁLÓ",
    2)

$rLabel = $d.Content
$found = $rLabel.Find.Execute(
    "FELHASZNÁLÓ",
    $true, $false, $false, $false, $false,
    $true, 1, $false)

if ($found) {
    $rColon = $d.Range($rLabel.End, $rLabel.End)
    $rColon.InsertAfter(":")
    $rColon.Bold = 1
}

# ---------------------------------------------------------------------
# Edit 3: append a new sentence about the "Admin" status to the player
# status paragraph.
# ---------------------------------------------------------------------
$rAdmin = $d.Content
$null = $rAdmin.Find.Execute(
    "kitiltás esetén “Kitiltott”.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "kitiltás esetén “Kitiltott”. Ha a státusz “Admin”, akkor a felhasználó kap szervezői jogokat.",
    2)
